$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q1" right after "2021-Q3" (and
#    before "总计"), populated with the new fund-holding data.
# ------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2021-Q3")
$newSheet = $wb.Worksheets.Add($null, $q3)
$newSheet.Name = "2022-Q1"

# Clone header formatting (bold / border / centered) from the 2021-Q3
# sheet so the new sheet reuses the identical style instead of
# creating a new one.
$q3.Range("B1:H1").Copy($newSheet.Range("B1:H1"))
$q3.Range("A2").Copy($newSheet.Range("A2"))

# Header row text (most labels match 2021-Q3 already; only the
# "基金规模" column name differs from "基金金额").
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data row (force text so numeric-looking values keep their original
# formatting / leading zeros instead of being coerced to numbers).
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").NumberFormat = "@"
$newSheet.Range("B2").Value = "005434"
$newSheet.Range("C2").Value = "鹏华睿投灵活配置混合"
$newSheet.Range("D2:G2").NumberFormat = "@"
$newSheet.Range("D2").Value = "3.41"
$newSheet.Range("E2").Value = "82.48"
$newSheet.Range("F2").Value = "1.76"
$newSheet.Range("G2").Value = "0.0600"
$newSheet.Range("H2").Value = 10

# ------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: add a new top row for
#    "2022-Q1" and push the existing rows down, renumbering the
#    index column (A) accordingly.
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Shift existing rows down (2020-Q4: row3->row4, 2021-Q3: row2->row3)
$total.Range("A3:D3").Copy($total.Range("A4:D4"))
$total.Range("A2:D2").Copy($total.Range("A3:D3"))

# Renumber the shifted rows' index column
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2

# New top row: 2022-Q1 summary
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.06
